$d = $word.ActiveDocument

# Locate the 4 target bullet paragraphs inside the "Experience" table cell by
# their (unique) current text.
function Find-ParaIndex($doc, $pattern) {
    $count = $doc.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        if ($doc.Paragraphs.Item($i).Range.Text -match $pattern) {
            return $i
        }
    }
    return -1
}

$idxResearched = Find-ParaIndex $d "^Researched object storage software called Ceph\."
if ($idxResearched -lt 0) { throw "Researched paragraph not found" }

# Step 1: insert one brand-new blank paragraph right after "Researched...".
# This gives us exactly 5 paragraph slots to hold the 5 final bullets, so the
# rest of the edit can proceed purely via single-paragraph replacements
# (no Cut/Delete/Move needed - those don't reliably persist in this host).
$pResearched = $d.Paragraphs.Item($idxResearched)
$pResearched.Range.InsertParagraphAfter()

$slot1 = $idxResearched        # -> "Researched ... Object Storage Software."
$slot2 = $idxResearched + 1    # -> "Built a Linux computer cluster ..."
$slot3 = $idxResearched + 2    # -> "Conducted read and write tests ..."
$slot4 = $idxResearched + 3    # -> "Contributed work to a team ..."
$slot5 = $idxResearched + 4    # -> "Presented findings ..."

$nsW = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$rPr = '<w:rPr><w:rFonts w:ascii="Nunito" w:hAnsi="Nunito"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr>'
$pPrTabs = '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:tabs><w:tab w:val="left" w:pos="2535"/></w:tabs>'
$pPrListStyle = '<w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>'

$xml1 = '<w:p ' + $nsW + '><w:pPr>' + $pPrTabs + $rPr + '</w:pPr>' + `
    '<w:r>' + $rPr + '<w:t xml:space="preserve">Researched </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/><w:r>' + $rPr + '<w:t>Ceph</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
    '<w:r>' + $rPr + '<w:t xml:space="preserve"> Object Storage Software.</w:t></w:r>' + `
    '</w:p>'

$xml2 = '<w:p ' + $nsW + '><w:pPr>' + $pPrListStyle + $rPr + '</w:pPr>' + `
    '<w:r>' + $rPr + '<w:t>Built a Linux computer cluster of 11 nodes and configured DNS, DHCP, NTP, NFS, HTTP, Authentication and Authorization, Users and Groups, LDAP, and Kerberos.</w:t></w:r>' + `
    '</w:p>'

$xml3 = '<w:p ' + $nsW + '><w:pPr>' + $pPrListStyle + $rPr + '</w:pPr>' + `
    '<w:r>' + $rPr + '<w:t xml:space="preserve">Conducted read and write tests on a computer cluster of 11 nodes. Increased performance of default configuration of </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/><w:r>' + $rPr + '<w:t>Ceph</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
    '<w:r>' + $rPr + '<w:t xml:space="preserve"> object storage software by 11%.</w:t></w:r>' + `
    '</w:p>'

$xml4 = '<w:p ' + $nsW + '><w:pPr>' + $pPrTabs + $rPr + '</w:pPr>' + `
    '<w:r>' + $rPr + '<w:t xml:space="preserve">Contributed work to a team of 3 computer science interns. </w:t></w:r>' + `
    '</w:p>'

$xml5 = '<w:p ' + $nsW + '><w:pPr>' + $pPrTabs + $rPr + '</w:pPr>' + `
    '<w:r>' + $rPr + '<w:t xml:space="preserve">Presented findings to Los Alamos National Laboratory computer scientists. </w:t></w:r>' + `
    '</w:p>'

# Replace each slot's whole paragraph (text + pPr) one at a time - InsertXML
# only behaves as a clean "replace" when the target Range is exactly one
# paragraph, so do these sequentially rather than as one multi-paragraph blob.
$d.Paragraphs.Item($slot1).Range.InsertXML($xml1)
$d.Paragraphs.Item($slot2).Range.InsertXML($xml2)
$d.Paragraphs.Item($slot3).Range.InsertXML($xml3)
$d.Paragraphs.Item($slot4).Range.InsertXML($xml4)
$d.Paragraphs.Item($slot5).Range.InsertXML($xml5)

Write-Output "done"
